$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.386.97"
$ws.Range("E2").Value = "  +7.53%  "
$ws.Range("D3").Value = "3.327.02"
$ws.Range("E3").Value = "  +2.83%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "409.89"
$ws.Range("E5").Value = "  +3.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.79"
$ws.Range("E6").Value = "  +6.01%  "
$ws.Range("D7").Value = "3.321.75"
$ws.Range("E7").Value = "  +2.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("E8").Value = "  -2.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.621"
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.111"
$ws.Range("E11").Value = "  +14.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "38.89"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").Value = "3.780.89"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.18"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.00"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "3.339.43"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").Value = "61.030.77"
$ws.Range("E18").Value = "  +7.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.991"
$ws.Range("E19").Value = "  -3.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.48"
$ws.Range("E20").Value = "  -3.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000113"
$ws.Range("E21").Value = "  +4.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.21"
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.40"
$ws.Range("E23").Value = "  -4.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "295.58"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.38"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.07"
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "28.71"
$ws.Range("E27").Value = "  +3.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.49"
$ws.Range("E28").Value = "  +3.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.40"
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.45"
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.109"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.15"
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "40.18"
$ws.Range("E35").Value = "  +4.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  +14.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0477"
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.50"
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.995"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("E40").Value = "  +4.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.31"
$ws.Range("E41").Value = "  -5.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "135.35"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.88"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.282"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.77"
$ws.Range("E46").Value = "  -4.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.16"
$ws.Range("E47").Value = "  -4.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.19"
$ws.Range("E48").Value = "  +3.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.95"
$ws.Range("E49").Value = "  -5.62%  "
$ws.Range("D50").Value = "2.114.68"
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("D51").Value = "3.647.53"
$ws.Range("E51").Value = "  +2.48%  "
